# ==========================================================================
# Slovenia Prva Liga - data refresh (21-04-2024 14:32)
#
# This script reproduces, via Excel COM automation, the same end result
# described by the authoritative OOXML diff:
#   1) Rows 9 and 10 (match ids 7 and 8, played on the same date) had their
#      entire data (every column except the running "id" in column A)
#      swapped between the two rows.
#   2) Two brand-new match rows were inserted right before the old row 154
#      (which shifts the former rows 154 and 155 down to 156 and 157),
#      and the new row 154 / row 155 were populated with the freshly
#      scraped match data.
# ==========================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------
# Step 1: swap the content of rows 9 and 10 (columns B through AC)
# --------------------------------------------------------------------
$row9 = @{
    "B"  = 6814328
    "F"  = "NK Domzale"
    "G"  = "NK Bravo"
    "I"  = 1
    "J"  = "D"
    "K"  = 2.35
    "L"  = 3.1
    "M"  = 2.9
    "N"  = 2.15
    "O"  = 3.1
    "P"  = 3.3
    "Q"  = -0.25
    "R"  = 1.925
    "S"  = 1.875
    "T"  = 2.25
    "U"  = 1.95
    "V"  = 1.85
    "W"  = -1
    "X"  = 2.1
    "AA" = 0.4375
    "AB" = -0.5
    "AC" = 0.425
}
foreach ($col in $row9.Keys) {
    $ws.Range(($col + "9")).Value2 = $row9[$col]
}

$row10 = @{
    "B"  = 6814330
    "F"  = "NK Maribor"
    "G"  = "NK Aluminij"
    "I"  = 0
    "J"  = "H"
    "K"  = 1.363
    "L"  = 4.5
    "M"  = 7
    "N"  = 1.4
    "O"  = 4.5
    "P"  = 7
    "Q"  = -1.25
    "R"  = 1.85
    "S"  = 1.95
    "T"  = 2.75
    "U"  = 1.8
    "V"  = 2
    "W"  = 0.3999999999999999
    "X"  = -1
    "AA" = 0.475
    "AB" = -1
    "AC" = 1
}
foreach ($col in $row10.Keys) {
    $ws.Range(($col + "10")).Value2 = $row10[$col]
}

# --------------------------------------------------------------------
# Step 2: insert two fresh rows right before the current row 154,
# pushing the old row 154 -> 156 and old row 155 -> 157 (their data
# travels with them automatically, untouched).
# --------------------------------------------------------------------
$ws.Range("154:155").Insert()

# Copy formatting (number format / font / borders) for the id (A) and
# Date (E) columns from the row directly below (the shifted old row
# 154, now row 156) so the two brand new rows look like every other
# data row in the table.
$ws.Range("A156").Copy()
$ws.Range("A154:A155").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E156").Copy()
$ws.Range("E154:E155").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# The running "id" sequence in column A must stay contiguous with the
# row position, so the two rows that shifted down (old 154 -> 156 and
# old 155 -> 157) need their id bumped by 2 as well.
$ws.Range("A156").Value2 = 154
$ws.Range("A157").Value2 = 155

# --------------------------------------------------------------------
# Step 3: populate the brand new row 154 (match id 152, Radomlje vs
# NS Mura, played 2024-04-21)
# --------------------------------------------------------------------
$ws.Range("A154").Value2 = 152
$ws.Range("B154").Value2 = 6994887
$ws.Range("C154").Value2 = "Slovenia Prva Liga"
$ws.Range("D154").Value2 = "Slovenia Prva Liga"
$ws.Range("E154").Value2 = 45402.41666666666
$ws.Range("F154").Value2 = "NK Radomlje"
$ws.Range("G154").Value2 = "NS Mura"
$ws.Range("H154").Value2 = 1
$ws.Range("I154").Value2 = 2
$ws.Range("J154").Value2 = "A"
$ws.Range("K154").Value2 = 2.625
$ws.Range("L154").Value2 = 3.25
$ws.Range("M154").Value2 = 2.45
$ws.Range("N154").Value2 = 2.1
$ws.Range("O154").Value2 = 3.25
$ws.Range("P154").Value2 = 3.1
$ws.Range("Q154").Value2 = -0.25
$ws.Range("R154").Value2 = 1.825
$ws.Range("S154").Value2 = 1.975
$ws.Range("T154").Value2 = 2.25
$ws.Range("U154").Value2 = 1.8
$ws.Range("V154").Value2 = 2
$ws.Range("W154").Value2 = -1
$ws.Range("X154").Value2 = -1
$ws.Range("Y154").Value2 = 2.1
$ws.Range("Z154").Value2 = -1
$ws.Range("AA154").Value2 = 0.9750000000000001
$ws.Range("AB154").Value2 = 0.8
$ws.Range("AC154").Value2 = -1

# --------------------------------------------------------------------
# Step 4: populate the brand new row 155 (match id 153, Olimpija
# Ljubljana vs NK Aluminij, played 2024-04-21)
# --------------------------------------------------------------------
$ws.Range("A155").Value2 = 153
$ws.Range("B155").Value2 = 6998172
$ws.Range("C155").Value2 = "Slovenia Prva Liga"
$ws.Range("D155").Value2 = "Slovenia Prva Liga"
$ws.Range("E155").Value2 = 45402.52083333334
$ws.Range("F155").Value2 = "Olimpija Ljubljana"
$ws.Range("G155").Value2 = "NK Aluminij"
$ws.Range("H155").Value2 = 5
$ws.Range("I155").Value2 = 0
$ws.Range("J155").Value2 = "H"
$ws.Range("K155").Value2 = 1.25
$ws.Range("L155").Value2 = 5.25
$ws.Range("M155").Value2 = 9.5
$ws.Range("N155").Value2 = 1.25
$ws.Range("O155").Value2 = 5.5
$ws.Range("P155").Value2 = 9
$ws.Range("Q155").Value2 = -1.75
$ws.Range("R155").Value2 = 1.875
$ws.Range("S155").Value2 = 1.925
$ws.Range("T155").Value2 = 3.25
$ws.Range("U155").Value2 = 1.95
$ws.Range("V155").Value2 = 1.85
$ws.Range("W155").Value2 = 0.25
$ws.Range("X155").Value2 = -1
$ws.Range("Y155").Value2 = -1
$ws.Range("Z155").Value2 = 0.875
$ws.Range("AA155").Value2 = -1
$ws.Range("AB155").Value2 = 0.95
$ws.Range("AC155").Value2 = -1
